$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for A2:B11 (header row A1:B1 "file 1"/"file 2" stays unchanged).
# Values are forced to Text (matching the original "numbers stored as text"
# shape of this sheet) via a leading apostrophe, then the cell style is put
# back to "Normal" so no stray number-format / quote-prefix style sticks to
# the cell (only the stored value's type changes).
$data = @(
    @("peter", "12"),
    @("2", "40"),
    @("3", "34"),
    @("4", "15"),
    @("5", "23"),
    @("6", "33"),
    @("7", "31"),
    @("8", "49"),
    @("9", "30"),
    @("10", "39")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2

    $cellA = $ws.Cells.Item($row, 1)
    $cellA.Value = "'" + $data[$i][0]
    $cellA.Style = "Normal"

    $cellB = $ws.Cells.Item($row, 2)
    $cellB.Value = "'" + $data[$i][1]
    $cellB.Style = "Normal"
}
